$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = 24.35904122085856
$ws.Range("D10").Value = 20.75097308944177
$ws.Range("D11").Value = 18.6108471309644
$ws.Range("D12").Value = 9.685567174640006
$ws.Range("D20").Value = 35.42061575668554
$ws.Range("D21").Value = 26.08107128583477
$ws.Range("D22").Value = 24.87166322780002
